$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a serial date value (45179 = 2023-09-10) for every
# data row (2-185). The update bumps this "changed" date forward by one day
# to 45180 (2023-09-11) for all of them.
$newValue = 45180

for ($row = 2; $row -le 185; $row++) {
    $ws.Cells.Item($row, 3).Value = $newValue
}
